# Updates the cryptos list (prices / 1h volume %) as scraped on
# Tue Jun 27 04:43:20 UTC 2023, applied via GitHub Actions.
#
# Columns: A=index, B=Coin, C=Link, D=Price, E=Volume(1h)
# D/E columns hold plain text (not real numbers) in the source workbook,
# so every write below forces the cell to Text format first and clears
# the style back to Normal afterwards so the cell keeps looking exactly
# like the untouched cells (no stray "@" number format / quote-prefix
# left behind, no float coercion / precision loss / scientific notation).

function Set-TextValue {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "30.380.71"
Set-TextValue $ws "E2" "  +0.33%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "1.874.10"
Set-TextValue $ws "E3" "  -0.23%  "

# Row 4 - TetherUSD
Set-TextValue $ws "E4" "  +0.00%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "238.53"

# Row 6 - USDC
Set-TextValue $ws "E6" "  +0.03%  "

# Row 7 - XRP
Set-TextValue $ws "D7" "0.4829"
Set-TextValue $ws "E7" "  -0.51%  "

# Row 8 - Cardano
Set-TextValue $ws "D8" "0.2826"
Set-TextValue $ws "E8" "  -1.61%  "

# Row 9 - Dogecoin
Set-TextValue $ws "D9" "0.06533"
Set-TextValue $ws "E9" "  -0.78%  "

# Row 10 - WrappedEther
Set-TextValue $ws "D10" "1.868.02"
Set-TextValue $ws "E10" "  -0.60%  "

# Row 11 - TRON
Set-TextValue $ws "D11" "0.07476"
Set-TextValue $ws "E11" "  +2.54%  "

# Row 12 - Solana
Set-TextValue $ws "D12" "16.54"
Set-TextValue $ws "E12" "  -1.17%  "

# Row 13 - Polkadot
Set-TextValue $ws "D13" "5.099"
Set-TextValue $ws "E13" "  -1.74%  "

# Row 14 - Litecoin
Set-TextValue $ws "D14" "88.18"
Set-TextValue $ws "E14" "  +1.23%  "

# Row 15 - Polygon
Set-TextValue $ws "D15" "0.6589"
Set-TextValue $ws "E15" "  +0.51%  "

# Row 16 - WrappedBTC
Set-TextValue $ws "D16" "30.357.05"
Set-TextValue $ws "E16" "  +0.37%  "

# Row 17 - Avalanche
Set-TextValue $ws "D17" "13.34"
Set-TextValue $ws "E17" "  -0.03%  "

# Row 18 - Dai
Set-TextValue $ws "E18" "  -0.03%  "

# Row 19 - ShibaInu
Set-TextValue $ws "D19" "0.000007636"
Set-TextValue $ws "E19" "  -1.17%  "

# Row 20 - WrappedliquidstakedEther2.0
Set-TextValue $ws "D20" "2.113.93"
Set-TextValue $ws "E20" "  -1.44%  "

# Row 21 - was BitcoinCash, now Uniswap (rows 21/22 swapped order)
Set-TextValue $ws "B21" "Uniswap"
Set-TextValue $ws "C21" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws "D21" "5.309"
Set-TextValue $ws "E21" "  +0.28%  "

# Row 22 - was Uniswap, now BitcoinCash
Set-TextValue $ws "B22" "BitcoinCash"
Set-TextValue $ws "C22" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws "D22" "222.89"
Set-TextValue $ws "E22" "  +14.38%  "

# Row 23 - BinanceUSD
Set-TextValue $ws "E23" "  +0.01%  "

# Row 24 - Chainlink
Set-TextValue $ws "D24" "6.199"
Set-TextValue $ws "E24" "  +1.14%  "

# Row 25 - Cosmos
Set-TextValue $ws "D25" "9.280"
Set-TextValue $ws "E25" "  -0.02%  "

# Row 26 - Monero
Set-TextValue $ws "D26" "165.93"
Set-TextValue $ws "E26" "  +4.31%  "

# Row 27 - EthereumClassic
Set-TextValue $ws "D27" "18.77"
Set-TextValue $ws "E27" "  +3.87%  "

# Row 28 - LidoDAOToken
Set-TextValue $ws "D28" "1.981"
Set-TextValue $ws "E28" "  +3.42%  "

# Row 29 - Toncoin
Set-TextValue $ws "D29" "1.459"
Set-TextValue $ws "E29" "  +1.33%  "

# Row 30 - Stellar
Set-TextValue $ws "D30" "0.09412"
Set-TextValue $ws "E30" "  +3.20%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue $ws "D31" "4.328"
Set-TextValue $ws "E31" "  +1.44%  "

# Row 32 - Filecoin
Set-TextValue $ws "D32" "4.023"
Set-TextValue $ws "E32" "  -0.75%  "

# Row 33 - Hedera
Set-TextValue $ws "D33" "0.05058"
Set-TextValue $ws "E33" "  -0.64%  "

# Row 34 - ARBITRUM
Set-TextValue $ws "D34" "1.218"
Set-TextValue $ws "E34" "  +11.32%  "

# Row 35 - ImmutableX
Set-TextValue $ws "D35" "0.7532"
Set-TextValue $ws "E35" "  +4.94%  "

# Row 36 - HuobiToken
Set-TextValue $ws "E36" "  -0.45%  "

# Row 37 - VeChain
Set-TextValue $ws "D37" "0.01843"
Set-TextValue $ws "E37" "  +2.64%  "

# Row 38 - MXToken
Set-TextValue $ws "D38" "2.619"
Set-TextValue $ws "E38" "  -0.62%  "

# Row 39 - RenderToken
Set-TextValue $ws "D39" "2.092"
Set-TextValue $ws "E39" "  +2.62%  "

# Row 40 - TrustWalletToken
Set-TextValue $ws "D40" "0.9075"
Set-TextValue $ws "E40" "  -1.32%  "

# Row 41 - FraxShare
Set-TextValue $ws "D41" "5.949"
Set-TextValue $ws "E41" "  +2.80%  "

# Row 42 - Quant
Set-TextValue $ws "D42" "107.00"
Set-TextValue $ws "E42" "  +0.66%  "

# Row 43 - TheSandbox
Set-TextValue $ws "D43" "0.4310"
Set-TextValue $ws "E43" "  +0.83%  "

# Row 44 - PaxDollar
Set-TextValue $ws "E44" "  +0.27%  "

# Row 45 - Aptos
Set-TextValue $ws "D45" "7.461"
Set-TextValue $ws "E45" "  +1.04%  "

# Row 46 - Algorand
Set-TextValue $ws "D46" "0.1303"
Set-TextValue $ws "E46" "  -1.23%  "

# Row 47 - Aave
Set-TextValue $ws "D47" "64.86"
Set-TextValue $ws "E47" "  -1.61%  "

# Row 48 - NEARProtocol
Set-TextValue $ws "D48" "1.497"
Set-TextValue $ws "E48" "  +9.73%  "

# Row 49 - EnergySwap
Set-TextValue $ws "D49" "9.076"
Set-TextValue $ws "E49" "  +1.65%  "

# Row 50 - Elrond
Set-TextValue $ws "D50" "34.23"
Set-TextValue $ws "E50" "  +1.15%  "

# Row 51 - was Cronos, now Decentraland
Set-TextValue $ws "B51" "Decentraland"
Set-TextValue $ws "C51" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws "D51" "0.3932"
Set-TextValue $ws "E51" "  +3.07%  "

Write-Output "cryptos list updated"
